$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B. This shifts the old B (annotation_name, wrap
# style) to C and the old C (annotation_type, always 0) to D, carrying their
# existing cell styles along (Excel's native column-insert behaviour), so no
# new style records are created.
$ws.Columns("B:B").Insert()

# The old "reference_id" column (D) is now at E after the insert above.
# Drop it entirely -- it is not part of the new layout. Deleting it shifts
# everything after back into place, so the blank formatted column F cells
# stay put (untouched) instead of sliding to G.
$ws.Columns("E:E").Delete()

# New column B header + values: segment_type (the coarse marker group that
# each fine-grained annotation_name belongs to).
$ws.Range("B1").Value = "segment_type"

$segmentType = @{
    2  = "PB2"
    3  = "PB1"
    4  = "PB1"
    5  = "PA"
    6  = "PA"
    7  = "HA"
    8  = "HA"
    9  = "NP"
    10 = "NA"
    11 = "MP"
    12 = "MP"
    13 = "NS"
    14 = "NS"
}

foreach ($row in $segmentType.Keys) {
    $ws.Cells.Item($row, 2).Value = $segmentType[$row]
}

# Match the author's final selection/cursor position.
$ws.Range("B14").Select()
